$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0133337242297736
$ws.Range("C2").Value = 0.9957916666666666
$ws.Range("D2").Value = 0.0516950312067346
$ws.Range("E2").Value = 0.9875833333333334
$ws.Range("F2").Value = 496.501354
$ws.Range("G2").Value = 87.256

$ws.Range("B3").Value = 0.013401185424634
$ws.Range("C3").Value = 0.9959895833333333
$ws.Range("D3").Value = 0.04504525618358698
$ws.Range("E3").Value = 0.9886250000000001
$ws.Range("F3").Value = 482.541859
$ws.Range("G3").Value = 87.28100000000001

$ws.Range("B4").Value = 0.01365568300724044
$ws.Range("C4").Value = 0.9958541666666667
$ws.Range("D4").Value = 0.04523311213725358
$ws.Range("E4").Value = 0.9881666666666667
$ws.Range("F4").Value = 480.79168
$ws.Range("G4").Value = 87.667

$ws.Range("B5").Value = 0.01328831710375107
$ws.Range("C5").Value = 0.9959791666666666
$ws.Range("D5").Value = 0.04574180978583525
$ws.Range("E5").Value = 0.9882916666666667
$ws.Range("F5").Value = 480.560282
$ws.Range("G5").Value = 87.529

$ws.Range("B6").Value = 0.01300265827286373
$ws.Range("C6").Value = 0.9960458333333333
$ws.Range("D6").Value = 0.04431977491540214
$ws.Range("E6").Value = 0.9889666666666667
$ws.Range("F6").Value = 479.49993
$ws.Range("G6").Value = 87.78700000000001

$ws.Range("B7").Value = 0.01321161453224037
$ws.Range("C7").Value = 0.9960069444444445
$ws.Range("D7").Value = 0.04389502741978821
$ws.Range("E7").Value = 0.9890694444444444
$ws.Range("F7").Value = 477.997565
$ws.Range("G7").Value = 88.021

$ws.Range("B8").Value = 0.01303147935201958
$ws.Range("C8").Value = 0.9960654761904763
$ws.Range("D8").Value = 0.04265030976048182
$ws.Range("E8").Value = 0.9894285714285714
$ws.Range("F8").Value = 479.871528
$ws.Range("G8").Value = 87.60899999999999

$ws.Range("B9").Value = 0.01309860517669904
$ws.Range("C9").Value = 0.9960390625
$ws.Range("D9").Value = 0.04444268784804419
$ws.Range("E9").Value = 0.9890833333333333
$ws.Range("F9").Value = 480.138173
$ws.Range("G9").Value = 87.214

$ws.Range("B10").Value = 0.01289298439493925
$ws.Range("C10").Value = 0.9961064814814815
$ws.Range("D10").Value = 0.0434578655803011
$ws.Range("E10").Value = 0.9894722222222223
$ws.Range("F10").Value = 480.097193
$ws.Range("G10").Value = 88.14400000000001

$ws.Range("B11").Value = 0.01298296638639392
$ws.Range("C11").Value = 0.9960770833333333
$ws.Range("D11").Value = 0.04330564696884187
$ws.Range("E11").Value = 0.9895416666666668
$ws.Range("F11").Value = 478.83235
$ws.Range("G11").Value = 87.485

